$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (ExpStartTime) changes for every data row (2-46) to the same new value
$ws.Range("D2:D46").Value = 32539.366537900001

# Row 2: StimOnsetTime / Response / RespTime updates, plus new Response cell
$ws.Range("F2").Value = 32544.446955200001
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.64923369999814895

# Row 3: StimOnsetTime / RespTime updates; Response cell removed
$ws.Range("F3").Value = 32552.920392100001
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 2.5004396999975143

# Row 4: StimOnsetTime / RespTime updates, plus new Response cell
$ws.Range("F4").Value = 32558.141877099999
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.36716630000228179

# Row 5: StimOnsetTime / Response / RespTime updates
$ws.Range("F5").Value = 32566.732453100001
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.88267259999702219

# Row 6: StimOnsetTime / RespTime updates; Response cell removed
$ws.Range("F6").Value = 32573.247048500001
$ws.Range("G6").ClearContents()
$ws.Range("H6").Value = 2.5005311999993864

Write-Output "done"
